$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New day-section header row (348): reuse the existing date-format style (s=1)
#     by copying formatting from a prior date cell, then write the actual values. ---
$ws.Cells.Item(40, 2).Copy()
$ws.Cells.Item(348, 2).PasteSpecial(-4122)

$ws.Range("B348").Value = 43496
$ws.Range("C348").Value = 1
$ws.Range("D348").Value = "momin"
$ws.Range("E348").Value = 415

$ws.Range("C349").Value = 2
$ws.Range("D349").Value = "farookh"
$ws.Range("E349").Value = 420

$ws.Range("C350").Value = 3
$ws.Range("D350").Value = "jabir"
$ws.Range("E350").Value = 419

$ws.Range("C351").Value = 4
$ws.Range("D351").Value = "hasen"
$ws.Range("E351").Value = 496

$ws.Range("C352").Value = 5
$ws.Range("D352").Value = "rasid"
$ws.Range("E352").Value = 410

$ws.Range("C353").Value = 6
$ws.Range("D353").Value = "khurseed"
$ws.Range("E353").Value = 455

$ws.Range("C354").Value = 7
$ws.Range("D354").Value = "toheed"
$ws.Range("E354").Value = 320

$ws.Range("C355").Value = 8
$ws.Range("D355").Value = "shokeen"
$ws.Range("E355").Value = 453

$ws.Range("C356").Value = 9
$ws.Range("D356").Value = "hasan"
$ws.Range("E356").Value = 538

$ws.Range("C357").Value = 10
$ws.Range("D357").Value = "isran"
$ws.Range("E357").Value = 470

$ws.Range("C358").Value = 11
$ws.Range("D358").Value = "shabir"
$ws.Range("E358").Value = 470

$ws.Range("C359").Value = 12
$ws.Range("D359").Value = "amjad"
$ws.Range("E359").Value = 411

$ws.Range("C360").Value = 13
$ws.Range("D360").Value = "wajid"
$ws.Range("E360").Value = 419

$ws.Range("C361").Value = 14
$ws.Range("D361").Value = "aashu"
$ws.Range("E361").Value = 430

$ws.Range("C362").Value = 15
$ws.Range("D362").Value = "aarif"
$ws.Range("E362").Value = 416

$ws.Range("C363").Value = 16
$ws.Range("D363").Value = "mustakim"
$ws.Range("E363").Value = 405

$ws.Range("C364").Value = 17
$ws.Range("D364").Value = "momin"
$ws.Range("E364").Value = 415

$ws.Range("C365").Value = 18
$ws.Range("D365").Value = "rafakat"
$ws.Range("E365").Value = 538

$ws.Range("C366").Value = 19
$ws.Range("D366").Value = "jabir"
$ws.Range("E366").Value = 419

$ws.Range("C367").Value = 20
$ws.Range("D367").Value = "shabir"
$ws.Range("E367").Value = 470

$ws.Range("C368").Value = 21
$ws.Range("D368").Value = "farookh"
$ws.Range("E368").Value = 420

$ws.Range("C369").Value = 22
$ws.Range("D369").Value = "rasid"
$ws.Range("E369").Value = 410

$ws.Range("C370").Value = 23
$ws.Range("D370").Value = "haseen"
$ws.Range("E370").Value = 496

$ws.Range("C371").Value = 24
$ws.Range("D371").Value = "shokeen"
$ws.Range("E371").Value = 453

$ws.Range("C372").Value = 25
$ws.Range("D372").Value = "irsad"
$ws.Range("E372").Value = 453

$ws.Range("C373").Value = 26
$ws.Range("D373").Value = "isran"
$ws.Range("E373").Value = 470

$ws.Range("C374").Value = 27
$ws.Range("D374").Value = "khurseed"
$ws.Range("E374").Value = 455

$ws.Range("C375").Value = 28
$ws.Range("D375").Value = "amjad"
$ws.Range("E375").Value = 411

$ws.Range("C376").Value = 29
$ws.Range("D376").Value = "aashu"
$ws.Range("E376").Value = 430

$ws.Range("C377").Value = 30
$ws.Range("D377").Value = "aarif"
$ws.Range("E377").Value = 416

$ws.Range("C378").Value = 31
$ws.Range("D378").Value = "hasan"
$ws.Range("E378").Value = 538

$ws.Range("C379").Value = 32
$ws.Range("D379").Value = "toheed"
$ws.Range("E379").Value = 320

$ws.Range("C380").Value = 33
$ws.Range("D380").Value = "momin"
$ws.Range("E380").Value = 415

$ws.Range("C381").Value = 34
$ws.Range("D381").Value = "jabir"
$ws.Range("E381").Value = 419

$ws.Range("C382").Value = 35
$ws.Range("D382").Value = "shabir"
$ws.Range("E382").Value = 470

$ws.Range("C383").Value = 36
$ws.Range("D383").Value = "mustakim"
$ws.Range("E383").Value = 405

$ws.Range("C384").Value = 37
$ws.Range("D384").Value = "rasid"
$ws.Range("E384").Value = 410

$ws.Range("C385").Value = 38
$ws.Range("D385").Value = "shokeen"
$ws.Range("E385").Value = 453

$ws.Range("C386").Value = 39
$ws.Range("D386").Value = "rafakat"
$ws.Range("E386").Value = 538

$ws.Range("C387").Value = 40
$ws.Range("D387").Value = "farookh"
$ws.Range("E387").Value = 420

$ws.Range("C388").Value = 41
$ws.Range("D388").Value = "haseen"
$ws.Range("E388").Value = 496

$ws.Range("C389").Value = 42
$ws.Range("D389").Value = "khurseed"
$ws.Range("E389").Value = 455

$ws.Range("C390").Value = 43
$ws.Range("D390").Value = "aashu"
$ws.Range("E390").Value = 430

$ws.Range("C391").Value = 44
$ws.Range("D391").Value = "irsad"
$ws.Range("E391").Value = 435

$ws.Range("C392").Value = 45
$ws.Range("D392").Value = "amjad"
$ws.Range("E392").Value = 411

$ws.Range("C393").Value = 46
$ws.Range("D393").Value = "isran"
$ws.Range("E393").Value = 470

$ws.Range("C394").Value = 47
$ws.Range("D394").Value = "toheed"
$ws.Range("E394").Value = 320

$ws.Range("C395").Value = 48
$ws.Range("D395").Value = "aarif"
$ws.Range("E395").Value = 416

$ws.Range("C396").Value = 49
$ws.Range("D396").Value = "hasan"
$ws.Range("E396").Value = 518

$ws.Range("C397").Value = 50
$ws.Range("D397").Value = "mehmood"
$ws.Range("E397").Value = 415

$ws.Range("C398").Value = 51
$ws.Range("D398").Value = "mehmood"
$ws.Range("E398").Value = 419

$ws.Range("C399").Value = 52
$ws.Range("D399").Value = "mehmood"
$ws.Range("E399").Value = 470

$ws.Range("C400").Value = 53
$ws.Range("D400").Value = "mehmood"
$ws.Range("E400").Value = 405

$ws.Range("E401").Formula = "=SUM(E2:E400)"

# --- Update the sheet view to match the scrolled/selected state after the append ---
$ws.Activate()
$excel.Goto($ws.Range("A378"), $true)
$excel.ActiveWindow.ScrollRow = 378
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D401").Select()
